# Apply the "Nb-rich" dataset additions extracted from 10.1016/j.matchar.2023.113301
# and a handful of small cleanups, per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Rows 169-179: fill in previously-blank rows with new data rows
#    for composition "(TiNbZr)29.67 (AlTa)5.5"
# -----------------------------------------------------------------

$composition = "(TiNbZr)29.67 (AlTa)5.5"
$structure   = "BCC+B2+Zr5Al3"
$processing  = "VAM+H+H"
$comment     = "homogenized at 750*C for 3 h and at 1050*C for 4 h"
$doi         = "10.1016/j.matchar.2023.113301"

function Fill-Row {
    param(
        [int]$r,
        [string]$name,
        [string]$source,
        [string]$params_,
        [double]$temp,
        [double]$value,
        $uncertainty,
        [string]$unit,
        [string]$pointer
    )
    $ws.Cells.Item($r, 2).Value = $composition
    $ws.Cells.Item($r, 3).Value = $structure
    $ws.Cells.Item($r, 4).Value = $processing
    $ws.Cells.Item($r, 5).Value = $comment
    $ws.Cells.Item($r, 6).Value = $name
    $ws.Cells.Item($r, 7).Value = $source
    $ws.Cells.Item($r, 8).Value = $params_
    $ws.Cells.Item($r, 9).Value = $temp
    $ws.Cells.Item($r, 10).Value = $value
    if ($null -ne $uncertainty) {
        $ws.Cells.Item($r, 11).Value = $uncertainty
    }
    $ws.Cells.Item($r, 12).Value = $unit
    if ($pointer -ne "") {
        $ws.Cells.Item($r, 13).Value = $pointer
    }
    $ws.Cells.Item($r, 14).Value = $doi
}

Fill-Row 169 "density"                      "EXP" "strain rate 1e-3/s" 298  6000       $null "kg/m^3" "T1"
Fill-Row 170 "UCS"                          "EXP" "strain rate 1e-3/s" 298  1025000000 $null "Pa"     ""
Fill-Row 171 "UCS"                          "EXP" "strain rate 1e-3/s" 873  670000000  $null "Pa"     ""
Fill-Row 172 "compressive yield stress"     "EXP" "strain rate 1e-3/s" 298  890000000  8000000 "Pa"   "T1"
Fill-Row 173 "compressive yield stress"     "EXP" "strain rate 1e-3/s" 873  625000000  5000000 "Pa"   "T1"
Fill-Row 174 "compressive yield stress"     "EXP" "strain rate 1e-3/s" 1073 210000000  8000000 "Pa"   "T1"
Fill-Row 175 "compressive yield stress"     "EXP" "strain rate 1e-3/s" 1273 60000000   10000000 "Pa"  "T1"
Fill-Row 176 "minimum compressive ductility" "EXP" "strain rate 1e-3/s" 298  70         $null "%"     ""
Fill-Row 177 "minimum compressive ductility" "EXP" "strain rate 1e-3/s" 873  70         $null "%"     ""
Fill-Row 178 "minimum compressive ductility" "EXP" "strain rate 1e-3/s" 1073 70         $null "%"     ""
Fill-Row 179 "minimum compressive ductility" "EXP" "strain rate 1e-3/s" 1273 70         $null "%"     ""

# -----------------------------------------------------------------
# 2) Row 258: remove stray composition/material-comment cells (B, E)
#    and normalize C/D/F/G/H to the plain "blank separator" style
#    used throughout the sheet (copy style from row 1's matching cols)
# -----------------------------------------------------------------

$ws.Cells.Item(258, 2).Clear()
$ws.Cells.Item(258, 5).Clear()
$ws.Cells.Item(1, 3).Copy($ws.Cells.Item(258, 3))
$ws.Cells.Item(1, 4).Copy($ws.Cells.Item(258, 4))
$ws.Cells.Item(1, 6).Copy($ws.Cells.Item(258, 6))
$ws.Cells.Item(1, 7).Copy($ws.Cells.Item(258, 7))
$ws.Cells.Item(1, 8).Copy($ws.Cells.Item(258, 8))

# -----------------------------------------------------------------
# 3) Row 305: A305 style changes from the "24" family to the "44"
#    family (matching the other blank-separator rows, e.g. row 429)
# -----------------------------------------------------------------

$ws.Cells.Item(429, 1).Copy($ws.Cells.Item(305, 1))

# -----------------------------------------------------------------
# 4) Remove trailing blank row 430 (dimension shrinks to A1:T429)
# -----------------------------------------------------------------

$ws.Rows.Item(430).Delete()

# -----------------------------------------------------------------
# 5) Update the view: scroll position and active cell selection
# -----------------------------------------------------------------

$excel.ActiveWindow.ScrollRow = 147
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H183").Select()
